# "took away parking lot filter"
#
# Remove the "Display Parking Lots" bullet and its "Press map marker for
# information on parking lots" sub-bullet from the Major Software Functions
# list. Everything below them (Display Bike Racks, Display Food Options,
# etc.) simply shifts up to fill the gap. The document's trailing "_GoBack"
# bookmark -- which Word drops at the point of the most recent edit -- moves
# from the end of the document to right before "Display Bike Racks" (the
# bullet that now sits where the deleted text used to be).

$d = $word.ActiveDocument

# Find the "Display Parking Lots" paragraph and the sub-bullet that follows it.
$headingPara = $null
$subPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.Trim()
    if ($t -eq "Display Parking Lots") {
        $headingPara = $p
    } elseif ($headingPara -ne $null -and $subPara -eq $null -and `
              $t -eq "Press map marker for information on parking lots") {
        $subPara = $p
    }
}

if ($headingPara -ne $null -and $subPara -ne $null) {
    # Delete both paragraphs (heading bullet + its sub-bullet) as one range.
    $deleteRange = $d.Range($headingPara.Range.Start, $subPara.Range.End)
    $deleteRange.Delete()
}

# Re-locate "Display Bike Racks" now that it has shifted up, and move the
# "_GoBack" bookmark to sit immediately before it.
$bikePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Display Bike Racks") {
        $bikePara = $p
        break
    }
}

if ($bikePara -ne $null) {
    $bookmarkRange = $d.Range($bikePara.Range.Start, $bikePara.Range.Start)
    $d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
}
